$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.01139766666666667
$ws.Cells.Item(2, 8).Value = 0.034193
$ws.Cells.Item(2, 9).Value = 0.1481191086775714
$ws.Cells.Item(2, 10).Value = 0.1481191086775714
$ws.Cells.Item(2, 13).Value = 0.344913
$ws.Cells.Item(2, 14).Value = 1.034739
$ws.Cells.Item(2, 15).Value = 0.07070596358783537
$ws.Cells.Item(2, 16).Value = 0.07070596358783539
$ws.Cells.Item(2, 17).Value = 0.003931203403
$ws.Cells.Item(2, 18).Value = 0.035380830627
$ws.Cells.Item(2, 19).Value = 0.01047290430481899
$ws.Cells.Item(2, 20).Value = 0.01047290430481899
$ws.Cells.Item(3, 7).Value = 0.01139766666666667
$ws.Cells.Item(3, 8).Value = 0.034193
$ws.Cells.Item(3, 9).Value = 0.1481191086775714
$ws.Cells.Item(3, 10).Value = 0.1481191086775714
$ws.Cells.Item(3, 15).Value = 0.1806419055929541
$ws.Cells.Item(3, 16).Value = 0.1806419055929541
$ws.Cells.Item(3, 17).Value = 0.01004356687833333
$ws.Cells.Item(3, 18).Value = 0.09039210190499999
$ws.Cells.Item(3, 19).Value = 0.02675651804624636
$ws.Cells.Item(3, 20).Value = 0.02675651804624637
$ws.Cells.Item(4, 7).Value = 0.01139766666666667
$ws.Cells.Item(4, 8).Value = 0.034193
$ws.Cells.Item(4, 9).Value = 0.1481191086775714
$ws.Cells.Item(4, 10).Value = 0.1481191086775714
$ws.Cells.Item(4, 13).Value = 2.052799333333333
$ws.Cells.Item(4, 14).Value = 6.158397999999999
$ws.Cells.Item(4, 15).Value = 0.420816712956019
$ws.Cells.Item(4, 16).Value = 0.420816712956019
$ws.Cells.Item(4, 17).Value = 0.02339712253488889
$ws.Cells.Item(4, 18).Value = 0.210574102814
$ws.Cells.Item(4, 19).Value = 0.06233099643967095
$ws.Cells.Item(4, 20).Value = 0.06233099643967095
$ws.Cells.Item(5, 7).Value = 0.01139766666666667
$ws.Cells.Item(5, 8).Value = 0.034193
$ws.Cells.Item(5, 9).Value = 0.1481191086775714
$ws.Cells.Item(5, 10).Value = 0.1481191086775714
$ws.Cells.Item(5, 13).Value = 0.2678313333333334
$ws.Cells.Item(5, 14).Value = 0.803494
$ws.Cells.Item(5, 15).Value = 0.05490449041453371
$ws.Cells.Item(5, 16).Value = 0.05490449041453372
$ws.Cells.Item(5, 17).Value = 0.003052652260222223
$ws.Cells.Item(5, 18).Value = 0.027473870342
$ws.Cells.Item(5, 19).Value = 0.008132404182596996
$ws.Cells.Item(5, 20).Value = 0.008132404182596996
$ws.Cells.Item(6, 7).Value = 0.01139766666666667
$ws.Cells.Item(6, 8).Value = 0.034193
$ws.Cells.Item(6, 9).Value = 0.1481191086775714
$ws.Cells.Item(6, 10).Value = 0.1481191086775714
$ws.Cells.Item(6, 13).Value = 1.331393
$ws.Cells.Item(6, 14).Value = 3.994179
$ws.Cells.Item(6, 15).Value = 0.2729309274486578
$ws.Cells.Item(6, 16).Value = 0.2729309274486578
$ws.Cells.Item(6, 17).Value = 0.01517477361633333
$ws.Cells.Item(6, 18).Value = 0.136572962547
$ws.Cells.Item(6, 19).Value = 0.04042628570423809
$ws.Cells.Item(6, 20).Value = 0.0404262857042381
$ws.Cells.Item(7, 7).Value = 0.015206
$ws.Cells.Item(7, 8).Value = 0.045618
$ws.Cells.Item(7, 9).Value = 0.197610548932631
$ws.Cells.Item(7, 10).Value = 0.197610548932631
$ws.Cells.Item(7, 13).Value = 0.344913
$ws.Cells.Item(7, 14).Value = 1.034739
$ws.Cells.Item(7, 15).Value = 0.07070596358783537
$ws.Cells.Item(7, 16).Value = 0.07070596358783539
$ws.Cells.Item(7, 17).Value = 0.005244747078
$ws.Cells.Item(7, 18).Value = 0.047202723702
$ws.Cells.Item(7, 19).Value = 0.01397224427740277
$ws.Cells.Item(7, 20).Value = 0.01397224427740277
$ws.Cells.Item(8, 7).Value = 0.015206
$ws.Cells.Item(8, 8).Value = 0.045618
$ws.Cells.Item(8, 9).Value = 0.197610548932631
$ws.Cells.Item(8, 10).Value = 0.197610548932631
$ws.Cells.Item(8, 15).Value = 0.1806419055929541
$ws.Cells.Item(8, 16).Value = 0.1806419055929541
$ws.Cells.Item(8, 17).Value = 0.01339945117
$ws.Cells.Item(8, 18).Value = 0.12059506053
$ws.Cells.Item(8, 19).Value = 0.03569674612446017
$ws.Cells.Item(8, 20).Value = 0.03569674612446017
$ws.Cells.Item(9, 7).Value = 0.015206
$ws.Cells.Item(9, 8).Value = 0.045618
$ws.Cells.Item(9, 9).Value = 0.197610548932631
$ws.Cells.Item(9, 10).Value = 0.197610548932631
$ws.Cells.Item(9, 13).Value = 2.052799333333333
$ws.Cells.Item(9, 14).Value = 6.158397999999999
$ws.Cells.Item(9, 15).Value = 0.420816712956019
$ws.Cells.Item(9, 16).Value = 0.420816712956019
$ws.Cells.Item(9, 17).Value = 0.03121486666266666
$ws.Cells.Item(9, 18).Value = 0.280933799964
$ws.Cells.Item(9, 19).Value = 0.08315782164726432
$ws.Cells.Item(9, 20).Value = 0.08315782164726433
$ws.Cells.Item(10, 7).Value = 0.015206
$ws.Cells.Item(10, 8).Value = 0.045618
$ws.Cells.Item(10, 9).Value = 0.197610548932631
$ws.Cells.Item(10, 10).Value = 0.197610548932631
$ws.Cells.Item(10, 13).Value = 0.2678313333333334
$ws.Cells.Item(10, 14).Value = 0.803494
$ws.Cells.Item(10, 15).Value = 0.05490449041453371
$ws.Cells.Item(10, 16).Value = 0.05490449041453372
$ws.Cells.Item(10, 17).Value = 0.004072643254666667
$ws.Cells.Item(10, 18).Value = 0.036653789292
$ws.Cells.Item(10, 19).Value = 0.01084970648968238
$ws.Cells.Item(10, 20).Value = 0.01084970648968239
$ws.Cells.Item(11, 7).Value = 0.015206
$ws.Cells.Item(11, 8).Value = 0.045618
$ws.Cells.Item(11, 9).Value = 0.197610548932631
$ws.Cells.Item(11, 10).Value = 0.197610548932631
$ws.Cells.Item(11, 13).Value = 1.331393
$ws.Cells.Item(11, 14).Value = 3.994179
$ws.Cells.Item(11, 15).Value = 0.2729309274486578
$ws.Cells.Item(11, 16).Value = 0.2729309274486578
$ws.Cells.Item(11, 17).Value = 0.020245161958
$ws.Cells.Item(11, 18).Value = 0.182206457622
$ws.Cells.Item(11, 19).Value = 0.05393403039382134
$ws.Cells.Item(11, 20).Value = 0.05393403039382136
$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 0.6666666666666666
$ws.Cells.Item(12, 7).Value = 0.05034566666666667
$ws.Cells.Item(12, 8).Value = 0.151037
$ws.Cells.Item(12, 9).Value = 0.6542703423897976
$ws.Cells.Item(12, 10).Value = 0.6542703423897976
$ws.Cells.Item(12, 13).Value = 0.344913
$ws.Cells.Item(12, 14).Value = 1.034739
$ws.Cells.Item(12, 15).Value = 0.07070596358783537
$ws.Cells.Item(12, 16).Value = 0.07070596358783539
$ws.Cells.Item(12, 17).Value = 0.017364874927
$ws.Cells.Item(12, 18).Value = 0.156283874343
$ws.Cells.Item(12, 19).Value = 0.04626081500561362
$ws.Cells.Item(12, 20).Value = 0.04626081500561362
$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = 0.6666666666666666
$ws.Cells.Item(13, 7).Value = 0.05034566666666667
$ws.Cells.Item(13, 8).Value = 0.151037
$ws.Cells.Item(13, 9).Value = 0.6542703423897976
$ws.Cells.Item(13, 10).Value = 0.6542703423897976
$ws.Cells.Item(13, 15).Value = 0.1806419055929541
$ws.Cells.Item(13, 16).Value = 0.1806419055929541
$ws.Cells.Item(13, 17).Value = 0.04436434973833334
$ws.Cells.Item(13, 18).Value = 0.399279147645
$ws.Cells.Item(13, 19).Value = 0.1181886414222476
$ws.Cells.Item(13, 20).Value = 0.1181886414222476
$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 0.6666666666666666
$ws.Cells.Item(14, 7).Value = 0.05034566666666667
$ws.Cells.Item(14, 8).Value = 0.151037
$ws.Cells.Item(14, 9).Value = 0.6542703423897976
$ws.Cells.Item(14, 10).Value = 0.6542703423897976
$ws.Cells.Item(14, 13).Value = 2.052799333333333
$ws.Cells.Item(14, 14).Value = 6.158397999999999
$ws.Cells.Item(14, 15).Value = 0.420816712956019
$ws.Cells.Item(14, 16).Value = 0.420816712956019
$ws.Cells.Item(14, 17).Value = 0.1033495509695556
$ws.Cells.Item(14, 18).Value = 0.9301459587259999
$ws.Cells.Item(14, 19).Value = 0.2753278948690838
$ws.Cells.Item(14, 20).Value = 0.2753278948690838
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 0.6666666666666666
$ws.Cells.Item(15, 7).Value = 0.05034566666666667
$ws.Cells.Item(15, 8).Value = 0.151037
$ws.Cells.Item(15, 9).Value = 0.6542703423897976
$ws.Cells.Item(15, 10).Value = 0.6542703423897976
$ws.Cells.Item(15, 13).Value = 0.2678313333333334
$ws.Cells.Item(15, 14).Value = 0.803494
$ws.Cells.Item(15, 15).Value = 0.05490449041453371
$ws.Cells.Item(15, 16).Value = 0.05490449041453372
$ws.Cells.Item(15, 17).Value = 0.01348414703088889
$ws.Cells.Item(15, 18).Value = 0.121357323278
$ws.Cells.Item(15, 19).Value = 0.03592237974225434
$ws.Cells.Item(15, 20).Value = 0.03592237974225434
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 0.6666666666666666
$ws.Cells.Item(16, 7).Value = 0.05034566666666667
$ws.Cells.Item(16, 8).Value = 0.151037
$ws.Cells.Item(16, 9).Value = 0.6542703423897976
$ws.Cells.Item(16, 10).Value = 0.6542703423897976
$ws.Cells.Item(16, 13).Value = 1.331393
$ws.Cells.Item(16, 14).Value = 3.994179
$ws.Cells.Item(16, 15).Value = 0.2729309274486578
$ws.Cells.Item(16, 16).Value = 0.2729309274486578
$ws.Cells.Item(16, 17).Value = 0.06702986818033334
$ws.Cells.Item(16, 18).Value = 0.603268813623
$ws.Cells.Item(16, 19).Value = 0.1785706113505983
$ws.Cells.Item(16, 20).Value = 0.1785706113505984
